$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Regenerate column G (K) values per recomputed strike counts (was using Strike#, now using K)
$ws.Range("G2").Value = 3
$ws.Range("G3").Value = 4
$ws.Range("G4").Value = 4
$ws.Range("G5").Value = 2
$ws.Range("G6").Value = 2
$ws.Range("G7").Value = 1
